$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the placeholder "(空)" text in H2 with a real hyperlink to the
# alum's homepage, matching the hyperlink style used elsewhere in the sheet.
$ws.Range("H2").Value = "https://tjuhaoxiaotian.github.io/"
$ws.Hyperlinks.Add($ws.Range("H2"), "https://tjuhaoxiaotian.github.io/") | Out-Null

# Match the look of the other hyperlink cells (underline + hyperlink theme color).
$ws.Range("H2").Font.Underline = 2
$ws.Range("H2").Font.Name = "Arial"
$ws.Range("H2").Font.Size = 10
$ws.Range("H2").Font.ThemeColor = 10

# Row 2 no longer needs the taller custom row height.
$ws.Rows.Item(2).RowHeight = $ws.Rows.Item(3).RowHeight

# Move the active selection, just reflecting where the author's cursor
# ended up after editing.
$ws.Range("H15").Select() | Out-Null
